# edit.ps1 - applies the "music -> biology" content rewrite described by the diff.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $find"
    }
}

# --- Title ---
Replace-Text "Harmonizing Nature's Melody: Music, Math, and the Aesthetics of Sound" "Unveiling the Symphony of Life: An Exploration of Biology"

# --- Author name / signature line ---
Replace-Text " Aurelio De Lira" " Sarah Morgan"

# --- Author first-name / e-mail line ---
Replace-Text "Aurelio" "sarah"
Replace-Text "DeLira@musicconservatory" "morgan@school"

# --- Body paragraph (intro) ---
Replace-Text "Within the vast tapestry of human existence, music occupies a realm of its own, weaving threads of emotion, intellect, and transcendence" "Biology, the study of life, is a captivating and complex field that encompasses the intricate workings of living organisms"
Replace-Text " Its allure lies not only in its ability to stir the soul but also in its profound connection to the intricate web of mathematics and the aesthetics of sound" " From the smallest microorganisms to the grandest whales, biology unveils the harmonious symphony of life and the profound interconnectedness of all living things"
Replace-Text " In this exploration, we delve into the harmonious interplay of music, math, and aesthetics, uncovering the hidden patterns that lend beauty and meaning to the melodies that grace our ears" " As we delve into this realm of knowledge, we embark on a journey of discovery, unraveling the mysteries of life and gaining a deeper understanding of ourselves and the natural world"

Replace-Text "From the subtle vibrations of a plucked string to the intricate harmonies of a symphony, music is governed by mathematical principles that create a foundation of order and symmetry" "Biology unveils the intricate mechanisms underlying the functioning of organisms, revealing the symphony between cells, tissues, organs, and systems"
Replace-Text " These mathematical underpinnings, such as ratios, proportions, and harmonic progressions, serve as the structural backbone of musical compositions, providing a framework for melodic development and variation" " Through meticulous observation and experimentation, biologists have unlocked the secrets of cellular respiration, DNA replication, and protein synthesis--the fundamental processes that sustain life"
Replace-Text " It is this mathematical order that enables musicians to craft cohesive pieces that resonate with our sense of equilibrium and balance" " These discoveries have revolutionized medicine, leading to innovative treatments and therapies"

Replace-Text "Moreover, music's aesthetic appeal stems from its ability to evoke emotions and convey narratives" "Biology not only probes the inner workings of organisms but also delves into the interactions between different species and the delicate balance of ecosystems"
Replace-Text " Whether it's the melancholic strains of a minor key or the uplifting melodies of a major chord, music possesses an uncanny ability to tap into our deepest emotions, triggering memories, and creating a sense of connection with others" " Studies of symbiotic relationships, food chains, and biogeochemical cycles underscore the interconnectedness of life and the importance of preserving biodiversity"

# This run keeps its own text; everything after it up through the "...hearts and minds..." run
# (several paragraphs' worth of "Introduction Continued/Concluded" material) is removed in the target,
# so first delete that trailing block, then rewrite the surviving run's text.
$text = $d.Content.Text

$lastKept = " This emotional resonance is further amplified by the aesthetic principles of form, texture, and dynamics, which composers employ to shape the overall structure and character of a piece, evoking a desired response from listeners"
$idxLastKept = $text.IndexOf($lastKept)
$endLastKept = $idxLastKept + $lastKept.Length

$lastRemoved = " It is through this interplay that music transcends its physical form, becoming a universal language capable of speaking to the hearts and minds of people across cultures, time, and space"
$idxLastRemoved = $text.IndexOf($lastRemoved)
$endLastRemoved = $idxLastRemoved + $lastRemoved.Length

$rng = $d.Range($endLastKept, $endLastRemoved)
$rng.Delete()

Replace-Text $lastKept " Biology equips us with the knowledge to address pressing environmental challenges and promote sustainable practices"

# --- Summary paragraph ---
Replace-Text "In this exploration of the relationship between music, mathematics, and aesthetics, we have uncovered the intricate web of principles that lend beauty and meaning to the melodies that grace our ears" "Biology is a captivating journey of discovery, unraveling the complexities of life and the interconnectedness of living organisms"
Replace-Text " Music's mathematical foundation provides a framework for creating cohesive compositions, while aesthetic principles shape its emotional impact and structure" " Through the study of cells, organisms, and ecosystems, we gain profound insights into the mechanisms underlying life and the intricate web of relationships that sustain our planet"
Replace-Text " This harmonious interplay enables music to transcend its physical form, becoming a universal language capable of speaking to the deepest recesses of the human experience" " Biology empowers us to address global challenges, understand our place in the natural world, and appreciate the mesmerizing symphony of life"

# --- Add a trailing empty paragraph after the Summary paragraph (before the section break) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

Write-Output "done"
